# add sprint 7 retrospective
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sprint 7 ("difficult to implement") retrospective counts bumped up
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 2

# Widen column A to fit the updated content, switching off the bestFit autosize
$ws.Columns("A").ColumnWidth = 41.8333333333333

# Update the view: zoom in and move the selection to C9
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("C9").Select()
